# Generate Report for Handback
#
# The localization-status report is being regenerated: the zh-cn and de-de
# handback for d3fcce66-c131-4993-82fc-f78e21534860 failed its transform
# (the handback file name didn't match the handoff file name), so:
#   - every "Ready for handoff" status cell for that file becomes
#     "Handback transform failed" (Overview!E3/F3, zh-cn!C3, de-de!C3 all
#     share that text)
#   - the per-locale "Error Detail" column (P) gets the failure detail text
#   - the Error Detail column is widened so the message is readable

$wb = $excel.ActiveWorkbook

$failedStatus = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $failedStatus
$wsOverview.Range("F3").Value = $failedStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $failedStatus
$wsZhCn.Range("P3").Value = "Handback file name: qfljuv2k.yec is different with handoff file name: d3fcce66-c131-4993-82fc-f78e21534860.d29761034705ea0e064b6ec38120a6377633f3e2.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $failedStatus
$wsDeDe.Range("P3").Value = "Handback file name: qfljuv2k.yec is different with handoff file name: d3fcce66-c131-4993-82fc-f78e21534860.d29761034705ea0e064b6ec38120a6377633f3e2.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
